$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-09-06 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-07 Saturday", 2)

$table = $d.Tables.Item(1)

$table.Cell(1,1).Range.Find.Execute("454÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "658÷7=", 2)
$table.Cell(1,2).Range.Find.Execute("582÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "291÷4=", 2)
$table.Cell(1,3).Range.Find.Execute("814÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "615÷9=", 2)
$table.Cell(1,4).Range.Find.Execute("314÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "295÷8=", 2)
$table.Cell(1,5).Range.Find.Execute("607÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "849÷4=", 2)

$table.Cell(5,1).Range.Find.Execute("476÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "945÷7=", 2)
$table.Cell(5,2).Range.Find.Execute("198÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "565÷5=", 2)
$table.Cell(5,3).Range.Find.Execute("598÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "817÷7=", 2)
$table.Cell(5,4).Range.Find.Execute("555÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "485÷2=", 2)
$table.Cell(5,5).Range.Find.Execute("140÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "501÷8=", 2)

$table.Cell(9,1).Range.Find.Execute("597÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "943÷4=", 2)
$table.Cell(9,2).Range.Find.Execute("672÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "912÷6=", 2)
$table.Cell(9,3).Range.Find.Execute("173÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "605÷7=", 2)
$table.Cell(9,4).Range.Find.Execute("826÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "811÷7=", 2)
$table.Cell(9,5).Range.Find.Execute("542÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "652÷3=", 2)

$table.Cell(13,1).Range.Find.Execute("830÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "351÷7=", 2)
$table.Cell(13,2).Range.Find.Execute("402÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "312÷9=", 2)
$table.Cell(13,3).Range.Find.Execute("596÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "542÷6=", 2)
$table.Cell(13,4).Range.Find.Execute("708÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "650÷9=", 2)
$table.Cell(13,5).Range.Find.Execute("876÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "337÷2=", 2)

$table.Cell(17,1).Range.Find.Execute("492÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "647÷9=", 2)
$table.Cell(17,2).Range.Find.Execute("823÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "542÷2=", 2)
$table.Cell(17,3).Range.Find.Execute("529÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "463÷5=", 2)
$table.Cell(17,4).Range.Find.Execute("705÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "200÷5=", 2)
$table.Cell(17,5).Range.Find.Execute("795÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "522÷7=", 2)
